$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.605527877807617
$ws.Range("B1").Value = 6.10035514831543
$ws.Range("C1").Value = 7.959798812866211
$ws.Range("D1").Value = 9.454241752624512
$ws.Range("E1").Value = 2.796086549758911
